$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D26:D182").Value = "Yes"
